# Fix the misspelled name "John Cenna" -> "John Cena" on the CarInsurance
# sheet, then leave that sheet active/selected (mirroring what Excel
# records when a user edits a cell on that sheet and saves).

$wb = $excel.ActiveWorkbook

$carInsurance = $wb.Worksheets.Item("CarInsurance")

# Correct the typo in the shared name value.
$carInsurance.Range("C2").Value = "John Cena"

# Make CarInsurance the active sheet, with G2 selected - reflects the
# cursor position left behind after editing the cell.
$carInsurance.Activate()
$carInsurance.Range("G2").Select()
